# Commit: "add some try catch when creating, reading file"
# Populate the three previously-blank "mynewsheet*" sheets with their text
# values, make "mynewsheet" the active/selected tab, and append three more
# new (still blank) worksheets: mynewsheet1231, Sheet5, Sheet6.

try {
    $wb = $excel.ActiveWorkbook

    # --- fill in the text values on the existing new sheets ---------------
    # Order matters: it drives the order entries land in sharedStrings.xml.
    try {
        $wsMynewsheet2 = $wb.Worksheets.Item("mynewsheet2")
        $wsMynewsheet2.Range("A1").Value = "sdfasf"
    } catch {
        Write-Host "failed writing mynewsheet2!A1: $_"
    }

    try {
        $wsMynewsheet3 = $wb.Worksheets.Item("mynewsheet3")
        $wsMynewsheet3.Range("A1").Value = "sassdasasdasd"
    } catch {
        Write-Host "failed writing mynewsheet3!A1: $_"
    }

    try {
        $wsMynewsheet = $wb.Worksheets.Item("mynewsheet")
        $wsMynewsheet.Range("A1").Value = "asdasd"
    } catch {
        Write-Host "failed writing mynewsheet!A1: $_"
    }

    # --- make "mynewsheet" the active tab ----------------------------------
    try {
        $wsMynewsheet.Activate()
    } catch {
        Write-Host "failed activating mynewsheet: $_"
    }

    # --- append three brand new (blank) worksheets at the end -------------
    $newSheetNames = @("mynewsheet1231", "Sheet5", "Sheet6")
    foreach ($newName in $newSheetNames) {
        try {
            $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
            $created = $wb.Worksheets.Add($null, $lastSheet)
            $created.Name = $newName
        } catch {
            Write-Host "failed creating sheet '$newName': $_"
        }
    }

    # re-activate mynewsheet since adding sheets moves selection to the
    # newest one
    try {
        $wsMynewsheet.Activate()
    } catch {
        Write-Host "failed re-activating mynewsheet: $_"
    }

    Write-Host "edit applied"
} catch {
    Write-Host "edit failed: $_"
}
